$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.211.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.317.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.79%  "
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.693.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.327.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.147.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.67%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -9.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.99%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.006.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.535.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("E51").Value = "  +1.19%  "
